# Apply the changes described by the diff:
#  1. Clear the empty placeholder cells B3 and B5 on the "ODI Batting" sheet
#     so they become truly blank (no cell record at all).
#  2. Add a new worksheet "ODI Batting Extra" (4th sheet, after "ODI Bowling")
#     with a header row and four data rows of stats.

$wb = $excel.ActiveWorkbook

# --- 1. Clear B3 / B5 on "ODI Batting" -------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("B3").ClearContents()
$battingSheet.Range("B5").ClearContents()

# --- 2. Add the new "ODI Batting Extra" sheet -------------------------------
# Worksheets.Add() inserts at the front by default, so pass the current
# last sheet as the "After" argument to place the new sheet at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Reuse the bold/centered/bordered header style already used by the other
# sheets (row 1, columns A:D on "Player Info") instead of constructing a
# brand-new style from scratch.
$headerSource = $wb.Worksheets.Item("Player Info").Range("A1:D1")
$headerSource.Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$newSheet.Cells.Item(1, 1).Value = "MATCH_CODE"
$newSheet.Cells.Item(1, 2).Value = "BATTING_POSITION"
$newSheet.Cells.Item(1, 3).Value = "NUM_4"
$newSheet.Cells.Item(1, 4).Value = "NUM_6"
$newSheet.Cells.Item(1, 5).Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Cells.Item(1, 6).Value = "MAN_OF_MATCH"

# Helper: write a value as genuine TEXT even when it looks like a number
# (Excel's COM layer auto-coerces plain numeric-looking strings assigned to
# .Value into real numbers; forcing the number format to Text ("@") first
# keeps it a string, matching the source inlineStr cells).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Row 2 - match 4597
Set-TextValue $newSheet.Cells.Item(2, 1) "4597"
$newSheet.Cells.Item(2, 2).Value = 10
Set-TextValue $newSheet.Cells.Item(2, 3) "0"
Set-TextValue $newSheet.Cells.Item(2, 4) "0"
Set-TextValue $newSheet.Cells.Item(2, 5) "0.53%"
$newSheet.Cells.Item(2, 6).Value = "NO"

# Row 3 - match 4600 (batting position / 4s / 6s / percent left blank)
Set-TextValue $newSheet.Cells.Item(3, 1) "4600"
$newSheet.Cells.Item(3, 6).Value = "NO"

# Row 4 - match 4601
Set-TextValue $newSheet.Cells.Item(4, 1) "4601"
$newSheet.Cells.Item(4, 2).Value = 10
Set-TextValue $newSheet.Cells.Item(4, 3) "3"
Set-TextValue $newSheet.Cells.Item(4, 4) "0"
Set-TextValue $newSheet.Cells.Item(4, 5) "5.91%"
$newSheet.Cells.Item(4, 6).Value = "NO"

# Row 5 - match 4603 (4s / 6s / percent left blank)
Set-TextValue $newSheet.Cells.Item(5, 1) "4603"
$newSheet.Cells.Item(5, 2).Value = 10
$newSheet.Cells.Item(5, 6).Value = "NO"
